$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (from 11.7109375 to 10.7109375).
# The engine quantizes ColumnWidth to steps of 1/6 character units, so
# 9.8 lands in the bucket that serializes to the closest achievable
# OOXML width (10.666666666666666) to the target 10.7109375.
$ws.Columns.Item(1).ColumnWidth = 9.8

# Update cell values
$ws.Range("A1").Value = 149.14734801281375
$ws.Range("B1").Value = 3.0534945669806746
$ws.Range("C1").Value = 3.9302022178734508
